# "remove fix bug from sorting"
# The "gsbmymohqweb09" row (row 7) in Sheet2 was a leftover duplicate/bug
# from a sort operation; delete it. Excel shifts the rows below up,
# re-adjusts the SUM formulas in the totals row, and drops the now-unused
# shared strings automatically.
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Rows.Item(7).Delete()

# Restore each sheet's view (zoom/selection) to match the saved workbook
# state, activating them in order so Sheet2 ends up the active tab.
$ws1 = $wb.Worksheets.Item("Sheet1")
[void]$ws1.Activate()
$excel.ActiveWindow.Zoom = 125
[void]$ws1.Range("A4:E4").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
[void]$ws3.Activate()
$excel.ActiveWindow.Zoom = 150
[void]$ws3.Range("A10:E10").Select()

$ws4 = $wb.Worksheets.Item("Sheet4")
[void]$ws4.Activate()
$excel.ActiveWindow.Zoom = 150
[void]$ws4.Range("A5:E5").Select()

[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 144
[void]$ws2.Range("A7:E16").Select()
